# Correction of ZP D and ZP S inversion: update the relative cover values
# for the "zp" (zoanthid plots) rows so that the Deep (D1-D3, rows 65-67)
# and Shallow (S1-S3, rows 68-70) records carry the correct algae,
# hard_coral, other, soft_coral, zoanthids, sd_alg and sd_hc values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 65 (zp_D1)
$ws.Range("F65").Value = 0.8
$ws.Range("H65").Value = 0.168
$ws.Range("L65").Value = 0.013
$ws.Range("M65").Value = 0.003
$ws.Range("N65").Value = 0.05
$ws.Range("P65").Value = 0.06

# Row 66 (zp_D2)
$ws.Range("F66").Value = 0.855
$ws.Range("H66").Value = 0.1
$ws.Range("L66").Value = 0.008999999999999999
$ws.Range("M66").Value = 0.001
$ws.Range("N66").Value = 0.05
$ws.Range("P66").Value = 0.06

# Row 67 (zp_D3)
$ws.Range("F67").Value = 0.901
$ws.Range("H67").Value = 0.047
$ws.Range("I67").Value = 0.004
$ws.Range("L67").Value = 0.024
$ws.Range("M67").Value = 0
$ws.Range("N67").Value = 0.05
$ws.Range("P67").Value = 0.06

# Row 68 (zp_S1)
$ws.Range("F68").Value = 0
$ws.Range("H68").Value = 0.86
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0.007
$ws.Range("N68").Value = 0.3
$ws.Range("P68").Value = 0.26

# Row 69 (zp_S2)
$ws.Range("F69").Value = 0.474
$ws.Range("H69").Value = 0.415
$ws.Range("L69").Value = 0.013
$ws.Range("M69").Value = 0.003
$ws.Range("N69").Value = 0.3
$ws.Range("P69").Value = 0.26

# Row 70 (zp_S3)
$ws.Range("F70").Value = 0.546
$ws.Range("H70").Value = 0.419
$ws.Range("I70").Value = 0.002
$ws.Range("L70").Value = 0.021
$ws.Range("M70").Value = 0.003
$ws.Range("N70").Value = 0.3
$ws.Range("P70").Value = 0.26
